# Correccion a Diebold Mariano y revision de Cap1
# Inserts the "d=6" column group (and its ARMA_I(*,6,*) rows) into the
# existing DeepAR results table, shifting d=7 (old col G) -> col H and
# d=10 (old col H) -> col I, and appending the corresponding new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Clone the existing header/label style (xf index 1: bold, boxed,
#        centered/top) onto the brand-new cells that need it, before any
#        values change. Using Copy + PasteSpecial(Formats) reuses the
#        existing style entry instead of minting new cellXfs rows. ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("A50").Copy() | Out-Null
$ws.Range("A51:A57").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2) Clear the data cells that become stale because their row's
#        label moved to a different column once the d=6 group was
#        inserted. ---
$staleCells = @(
  "H5","H6","H7",
  "B23","B24","H26",
  "C27","C28","D29","D30","E31","E32","F33","F34","G35","G36",
  "B37","B38","H39",
  "C41","C42","D43","D44","E45","E46","F47","F48","G49","G50"
)
foreach ($ref in $staleCells) {
  $ws.Range($ref).ClearContents()
}

# --- 3) Header row: shift d=7/d=10 right one column and add d=6. ---
$ws.Range("G1").Value = "d=6"
$ws.Range("H1").Value = "d=7"
$ws.Range("I1").Value = "d=10"

# --- 4) d=10 column values shift from H to I (rows 5-7, 28-29, 44-45). ---
$ws.Range("I5").Value = 96.00065223383682
$ws.Range("I6").Value = 95.82045501397208
$ws.Range("I7").Value = 95.88891128584875

# --- 5) New ARMA_I(0,6,*) rows (20-22) carrying the d=6 values that used
#        to sit in G20:G22 for ARMA_I(0,7,*). ---
$ws.Range("A20").Value = "ARMA_I(0,6,0)"
$ws.Range("G20").Value = 97.81293379130943
$ws.Range("A21").Value = "ARMA_I(0,6,1)"
$ws.Range("G21").Value = 97.9659223869358
$ws.Range("A22").Value = "ARMA_I(0,6,2)"
$ws.Range("G22").Value = 97.82715320351596

# --- 6) ARMA_I(0,7,*) rows, now at 23-25, values moved from G to H. ---
$ws.Range("A23").Value = "ARMA_I(0,7,0)"
$ws.Range("H23").Value = 97.28172209682585
$ws.Range("A24").Value = "ARMA_I(0,7,1)"
$ws.Range("H24").Value = 97.4173951488175
$ws.Range("A25").Value = "ARMA_I(0,7,2)"
$ws.Range("H25").Value = 97.26205355354016

# --- 7) ARMA_I(1,1,*) rows shift down to 26-27. ---
$ws.Range("A26").Value = "ARMA_I(1,1,0)"
$ws.Range("B26").Value = 59.60157916411891
$ws.Range("A27").Value = "ARMA_I(1,1,1)"
$ws.Range("B27").Value = 57.55437043155216

# --- 8) ARMA_I(1,10,*) rows shift to 28-29, values move from H to I. ---
$ws.Range("A28").Value = "ARMA_I(1,10,0)"
$ws.Range("I28").Value = 95.95145035243536
$ws.Range("A29").Value = "ARMA_I(1,10,1)"
$ws.Range("I29").Value = 95.86430176628636

# --- 9) ARMA_I(1,2,*) rows shift to 30-31. ---
$ws.Range("A30").Value = "ARMA_I(1,2,0)"
$ws.Range("C30").Value = 97.13262743084054
$ws.Range("A31").Value = "ARMA_I(1,2,1)"
$ws.Range("C31").Value = 96.98581963689303

# --- 10) ARMA_I(1,3,*) rows shift to 32-33. ---
$ws.Range("A32").Value = "ARMA_I(1,3,0)"
$ws.Range("D32").Value = 99.43323346220409
$ws.Range("A33").Value = "ARMA_I(1,3,1)"
$ws.Range("D33").Value = 99.39900532546699

# --- 11) ARMA_I(1,4,*) rows shift to 34-35. ---
$ws.Range("A34").Value = "ARMA_I(1,4,0)"
$ws.Range("E34").Value = 98.86214405837543
$ws.Range("A35").Value = "ARMA_I(1,4,1)"
$ws.Range("E35").Value = 98.81102535160306

# --- 12) ARMA_I(1,5,*) rows shift to 36-37. ---
$ws.Range("A36").Value = "ARMA_I(1,5,0)"
$ws.Range("F36").Value = 98.32671934964857
$ws.Range("A37").Value = "ARMA_I(1,5,1)"
$ws.Range("F37").Value = 98.52472330275019

# --- 13) New ARMA_I(1,6,*) rows (38-39). ---
$ws.Range("A38").Value = "ARMA_I(1,6,0)"
$ws.Range("G38").Value = 97.86057614725064
$ws.Range("A39").Value = "ARMA_I(1,6,1)"
$ws.Range("G39").Value = 97.8107354485456

# --- 14) ARMA_I(1,7,*) rows shift to 40-41, values move from G to H. ---
$ws.Range("A40").Value = "ARMA_I(1,7,0)"
$ws.Range("H40").Value = 97.22769138562258
$ws.Range("A41").Value = "ARMA_I(1,7,1)"
$ws.Range("H41").Value = 97.28389908304054

# --- 15) ARMA_I(2,1,*) rows shift to 42-43. ---
$ws.Range("A42").Value = "ARMA_I(2,1,0)"
$ws.Range("B42").Value = 71.19128083268376
$ws.Range("A43").Value = "ARMA_I(2,1,2)"
$ws.Range("B43").Value = 92.82561145747597

# --- 16) ARMA_I(2,10,*) rows shift to 44-45, values move from H to I. ---
$ws.Range("A44").Value = "ARMA_I(2,10,0)"
$ws.Range("I44").Value = 95.94094814452673
$ws.Range("A45").Value = "ARMA_I(2,10,2)"
$ws.Range("I45").Value = 95.85318798002508

# --- 17) ARMA_I(2,2,*) rows shift to 46-47. ---
$ws.Range("A46").Value = "ARMA_I(2,2,0)"
$ws.Range("C46").Value = 97.80091765826667
$ws.Range("A47").Value = "ARMA_I(2,2,2)"
$ws.Range("C47").Value = 96.81630305973297

# --- 18) ARMA_I(2,3,*) rows shift to 48-49. ---
$ws.Range("A48").Value = "ARMA_I(2,3,0)"
$ws.Range("D48").Value = 99.41623731688573
$ws.Range("A49").Value = "ARMA_I(2,3,2)"
$ws.Range("D49").Value = 99.40704935094892

# --- 19) ARMA_I(2,4,*) rows shift to 50-51 (row 51 is brand new). ---
$ws.Range("A50").Value = "ARMA_I(2,4,0)"
$ws.Range("E50").Value = 98.85595082336374
$ws.Range("A51").Value = "ARMA_I(2,4,2)"
$ws.Range("E51").Value = 98.84051625468612

# --- 20) ARMA_I(2,5,*) rows -> 52-53 (new rows). ---
$ws.Range("A52").Value = "ARMA_I(2,5,0)"
$ws.Range("F52").Value = 98.34031324133727
$ws.Range("A53").Value = "ARMA_I(2,5,2)"
$ws.Range("F53").Value = 98.215466319863

# --- 21) New ARMA_I(2,6,*) rows -> 54-55 (new rows). ---
$ws.Range("A54").Value = "ARMA_I(2,6,0)"
$ws.Range("G54").Value = 97.81815998633738
$ws.Range("A55").Value = "ARMA_I(2,6,2)"
$ws.Range("G55").Value = 97.80418473467401

# --- 22) ARMA_I(2,7,*) rows -> 56-57 (new rows), values move from G to H. ---
$ws.Range("A56").Value = "ARMA_I(2,7,0)"
$ws.Range("H56").Value = 97.32984359253554
$ws.Range("A57").Value = "ARMA_I(2,7,2)"
$ws.Range("H57").Value = 97.40925303193708
